$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price (col D), new Volume(1h) (col E), new Hora (col G).
# $null means that column is unchanged for that row (matches the source diff).
$updates = @(
    @(2, '317.17', '4.61%', '14'),
    @(3, '48.07', '11.57%', '14'),
    @(4, '5.251', '4.03%', '14'),
    @(5, '0.07981', '3.93%', '14'),
    @(6, '4.595', '4.24%', '14'),
    @(7, '1.429', '35.26%', '14'),
    @(8, '1.645', '2.11%', '14'),
    @(9, '0.1278', '3.59%', '14'),
    @(10, '0.1934', '4.44%', '14'),
    @(11, '0.09277', '2.39%', '14'),
    @(12, '0.04580', '10.03%', '14'),
    @(13, '0.1043', '-0.37%', '14'),
    @(14, '0.001320', '3.15%', '14'),
    @(15, '0.04171', '0.71%', '14'),
    @(16, '0.005925', '2.68%', '14'),
    @(17, '3.330', '0.09%', '14'),
    @(18, '2.432', '2.03%', '14'),
    @(19, '0.3470', '3.68%', '14'),
    @(20, '8.121', '-3.16%', '14'),
    @(21, '0.1404', '-0.13%', '14'),
    @(22, '0.3105', '7.40%', '14'),
    @(23, '0.001311', '2.77%', '14'),
    @(24, '0.004229', '-5.83%', '14'),
    @(25, '0.0001352', '0.31%', '14'),
    @(26, '0.0003546', $null, '14'),
    @(27, $null, $null, '14'),
    @(28, $null, $null, '14'),
    @(29, $null, $null, '14'),
    @(30, $null, $null, '14'),
    @(31, $null, $null, '14'),
    @(32, $null, $null, '14'),
    @(33, $null, $null, '14'),
    @(34, $null, $null, '14'),
    @(35, $null, $null, '14'),
    @(36, $null, $null, '14'),
    @(37, $null, $null, '14'),
    @(38, '0.02669', '8.81%', '14'),
    @(39, '0.05639', '6.86%', '14'),
    @(40, '0.008215', '36.82%', '14'),
    @(41, '0.008147', '6.42%', '14'),
    @(42, '0.1435', '6.35%', '14'),
    @(43, '0.007694', '4.69%', '14'),
    @(44, '0.008497', '14.39%', '14'),
    @(45, '0.3468', '14.71%', '14'),
    @(46, '0.00006919', '4.00%', '14'),
    @(47, '0.00000000751', '0.43%', '14'),
    @(48, '0.05496', '43.04%', '14'),
    @(49, '0.004008', '-4.60%', '14'),
    @(50, '0.00002104', '0.43%', '14'),
    @(51, '0.0002004', '0.43%', '14')
)

foreach ($u in $updates) {
    $row = $u[0]
    $newPrice  = $u[1]
    $newVolume = $u[2]
    $newHora   = $u[3]

    # Leading "'" keeps the value a literal text string (matches the workbook's
    # existing inline-string cells) instead of letting Excel coerce it to a number/percentage.
    if ($null -ne $newPrice)  { $ws.Cells.Item($row, 4).Value = "'" + $newPrice }
    if ($null -ne $newVolume) { $ws.Cells.Item($row, 5).Value = "'" + $newVolume }
    if ($null -ne $newHora)   { $ws.Cells.Item($row, 7).Value = "'" + $newHora }
}
